# Implement state specific work multipliers (#20)
#
# Translate the German federal-state ("Bundesland") names in column A to
# their English equivalents, and restore the last-used selection on the
# frozen (right) pane.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# German -> English state-name translations.
$translations = @{
    "Bayern"              = "Bavaria"
    "Niedersachsen"       = "Lower Saxony"
    "Nordrhein-Westfalen" = "North Rhine-Westphalia"
    "Rheinland-Pfalz"     = "Rhineland-Palatinate"
    "Sachsen"             = "Saxony"
    "Sachsen-Anhalt"      = "Saxony-Anhalt"
    "Thüringen"           = "Thuringia"
}

# Column A holds the "Bundesland" values (row 1 is the header).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Text
    if ($translations.ContainsKey($current)) {
        $cell.Value = $translations[$current]
    }
}

# Restore the active cell on the frozen right pane to H20 (was K17).
$ws.Range("H20").Select()
